$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.02110507534020301
$ws.Range("E4").Value = 0.01299622419984897
$ws.Range("F4").Value = 0.01845261203410448
$ws.Range("H4").Value = -0.01932280886891235
$ws.Range("J4").Value = -0.0002583969705665081

$ws.Range("C5").Value = -0.0103843876153755
$ws.Range("E5").Value = 0.0008820878752835149
$ws.Range("F5").Value = -0.0008629962585198503
$ws.Range("H5").Value = 0.001044197129767885
$ws.Range("J5").Value = 0.004634320760550871

$ws.Range("C6").Value = 0.005599665823986632
$ws.Range("E6").Value = -0.01079733825589353
$ws.Range("F6").Value = 0.00434991626999665
$ws.Range("H6").Value = -0.004338794477551778
$ws.Range("J6").Value = 0.006102731144028873

$ws.Range("C7").Value = 0.003708964084358563
$ws.Range("E7").Value = -0.006502285796091431
$ws.Range("F7").Value = -0.005087549291501971
$ws.Range("H7").Value = 0.005263921650556866
$ws.Range("J7").Value = 0.0178803832675434

$ws.Range("C8").Value = 0.1003132927485317
$ws.Range("E8").Value = -0.1291804318872173
$ws.Range("F8").Value = -0.9990109868404393
$ws.Range("H8").Value = 0.9999999263679969
$ws.Range("J8").Value = -0.01327333373497799

$ws.Range("C9").Value = 0.9556117589924702
$ws.Range("E9").Value = 0.009484258075370322
$ws.Range("F9").Value = 0.02704918293796731
$ws.Range("H9").Value = -0.0272487874899515
$ws.Range("J9").Value = 0.01514218753701872

$ws.Range("C10").Value = 0.00420973351238934
$ws.Range("E10").Value = -0.0256102650884106
$ws.Range("F10").Value = -0.0106556674022267
$ws.Range("H10").Value = 0.01121277203251088
$ws.Range("J10").Value = 0.004729352569465431

$ws.Range("C11").Value = -0.002767640078705602
$ws.Range("E11").Value = 0.01176585791063432
$ws.Range("F11").Value = -0.005739736453589458
$ws.Range("H11").Value = 0.005150287886011515
$ws.Range("J11").Value = 0.00938861510029193

$ws.Range("C12").Value = 0.06583375242535008
$ws.Range("E12").Value = 0.005407630680305227
$ws.Range("F12").Value = -0.02376814040672561
$ws.Range("H12").Value = 0.02384326789773071
$ws.Range("J12").Value = 0.02053901750261387

$ws.Range("C13").Value = 0.1040558466582338
$ws.Range("E13").Value = -0.007785995351439813
$ws.Range("F13").Value = 0.01314400084576003
$ws.Range("H13").Value = -0.01319955268798211
$ws.Range("J13").Value = 0.003063219359999851

$ws.Range("C14").Value = -0.2009269572210783
$ws.Range("E14").Value = -0.006815835248633409
$ws.Range("F14").Value = -0.02043957432158297
$ws.Range("H14").Value = 0.0207165825246633
$ws.Range("J14").Value = -0.0103139570346046

$ws.Range("C15").Value = -0.01358543574341743
$ws.Range("E15").Value = 0.007478063723122547
$ws.Range("F15").Value = -0.003017578104703124
$ws.Range("H15").Value = 0.002034005841360233
$ws.Range("J15").Value = -0.0226197537969832

$ws.Range("C16").Value = -0.005735100517404021
$ws.Range("E16").Value = 0.02060028216801129
$ws.Range("F16").Value = -0.02791851894074076
$ws.Range("H16").Value = 0.02683653131346125
$ws.Range("J16").Value = 0.01892102930053914

$ws.Range("C17").Value = 0.008202406696096266
$ws.Range("E17").Value = -0.01957260606290424
$ws.Range("F17").Value = -0.04102319741692789
$ws.Range("H17").Value = 0.04119831812793272
$ws.Range("J17").Value = -0.01581977773879409

$ws.Range("C18").Value = 0.02768263992330559
$ws.Range("E18").Value = -0.007037450393498015
$ws.Range("F18").Value = -0.0136662976026519
$ws.Range("H18").Value = 0.01336968072678723
$ws.Range("J18").Value = 0.0037378436055706

$ws.Range("C19").Value = 0.0113900300876012
$ws.Range("E19").Value = 0.008551873398074934
$ws.Range("F19").Value = -0.002496253155850126
$ws.Range("H19").Value = 0.002057380306295212
$ws.Range("J19").Value = 0.007263896298911399

$ws.Range("C20").Value = 0.009783588487343538
$ws.Range("E20").Value = 0.02644106198564248
$ws.Range("F20").Value = -0.000173226726929069
$ws.Range("H20").Value = -0.0005817705832708231
$ws.Range("J20").Value = -0.005287484288086148

$ws.Range("C21").Value = 0.02397818927912757
$ws.Range("E21").Value = -0.02289423835576954
$ws.Range("F21").Value = -0.0238473775618951
$ws.Range("H21").Value = 0.02410961184438447
$ws.Range("J21").Value = 0.006818723469679011

$ws.Range("C22").Value = 0.01142120330484813
$ws.Range("E22").Value = 0.01957146721485869
$ws.Range("F22").Value = 0.0004406947376277894
$ws.Range("H22").Value = -0.000767975550719022
$ws.Range("J22").Value = 0.004971578136974088

$ws.Range("C23").Value = -0.008519180980767238
$ws.Range("E23").Value = -0.001915117228604689
$ws.Range("F23").Value = 0.009148524845940992
$ws.Range("H23").Value = -0.00849640152385606
$ws.Range("J23").Value = 0.003500357872446744
